# Add new columns I (I0) and J (IF) to the sheet, mirroring the style of
# the existing header row / data cells, and populate their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — same bordered/bold style (StyleIndex 1) as the
# other header cells such as H1. Copy formats from H1 so we reuse the
# existing style entry rather than creating a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-21.
$values = @{
    2  = @(4, 5)
    3  = @(7, 7)
    4  = @(1, 4)
    5  = @(1, 6)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 3)
    9  = @(1, 6)
    10 = @(1, 7)
    11 = @(1, 6)
    12 = @(6, 8)
    13 = @(1, 5)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 6)
    18 = @(6, 8)
    19 = @(4, 7)
    20 = @(1, 3)
    21 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
